$d = $word.ActiveDocument
$replacements = @(
    @("2025-12-20 Saturday", "2025-12-21 Sunday"),
    @("12-2=10", "76-65=11"),
    @("30+65=95", "57-3=54"),
    @("82-63=19", "66-57=9"),
    @("80-52=28", "20+25=45"),
    @("80-36=44", "36+9=45"),
    @("25+62=87", "21-2=19"),
    @("50+49=99", "43+7=50"),
    @("88-61=27", "25+5=30"),
    @("0+33=33", "76-60=16"),
    @("58-52=6", "62+15=77"),
    @("94-56=38", "21+55=76"),
    @("88-58=30", "89-60=29"),
    @("80-76=4", "27+45=72"),
    @("81-2=79", "7+5=12"),
    @("51+4=55", "31+21=52"),
    @("28+11=39", "87+7=94"),
    @("42-8=34", "91-56=35"),
    @("26+0=26", "85-56=29"),
    @("70+26=96", "77+2=79"),
    @("57+9=66", "42-2=40"),
    @("1+44=45", "89-38=51"),
    @("33-1=32", "99-76=23"),
    @("47-14=33", "27-16=11"),
    @("84-44=40", "78-29=49"),
    @("30+56=86", "75-64=11"),
    @("13+10=23", "76-62=14"),
    @("48-28=20", "94-83=11"),
    @("55+14=69", "71-54=17"),
    @("86-79=7", "98-64=34"),
    @("30-0=30", "40+40=80"),
    @("98-40=58", "95-25=70"),
    @("32+22=54", "33+18=51"),
    @("56-53=3", "57+41=98"),
    @("30-22=8", "82-70=12"),
    @("39+47=86", "40+6=46"),
    @("62-60=2", "34-25=9"),
    @("0+60=60", "57-4=53"),
    @("85-65=20", "85-51=34"),
    @("10+81=91", "57+5=62"),
    @("6+70=76", "22-7=15"),
    @("12+26=38", "33+18=51"),
    @("13+45=58", "47-26=21"),
    @("22+69=91", "92-25=67"),
    @("36-35=1", "90-69=21"),
    @("1+84=85", "52+19=71"),
    @("29+50=79", "0+84=84"),
    @("4+39=43", "48+8=56"),
    @("19+9=28", "60+23=83"),
    @("28-11=17", "43-23=20"),
    @("31+20=51", "57-31=26"),
    @("51+23=74", "6+63=69"),
    @("2+34=36", "8+37=45"),
    @("74-0=74", "96-23=73"),
    @("72+23=95", "30+41=71"),
    @("52+20=72", "77-49=28"),
    @("80-78=2", "89-78=11"),
    @("61+7=68", "12+1=13"),
    @("46+8=54", "77-60=17"),
    @("93-33=60", "49+1=50"),
    @("45-1=44", "41+21=62"),
    @("16+23=39", "58-3=55"),
    @("56-45=11", "86-50=36"),
    @("61+34=95", "30+25=55"),
    @("43-3=40", "8+59=67"),
    @("31+5=36", "65-14=51"),
    @("35+40=75", "60+22=82"),
    @("40-10=30", "24+41=65"),
    @("55+13=68", "95-27=68"),
    @("73-47=26", "8+22=30"),
    @("24-1=23", "17+61=78"),
    @("77-52=25", "83-47=36"),
    @("76-16=60", "94-1=93"),
    @("25+1=26", "89-23=66"),
    @("61+4=65", "58-7=51"),
    @("71-52=19", "80-46=34"),
    @("49-9=40", "81-57=24"),
    @("39-33=6", "4+94=98"),
    @("39-8=31", "64+28=92"),
    @("25+18=43", "38+46=84"),
    @("74-15=59", "55-2=53"),
    @("6-3=3", "18+31=49"),
    @("11+1=12", "68-26=42"),
    @("9+24=33", "38-36=2"),
    @("2+40=42", "2+44=46"),
    @("73-15=58", "8+6=14"),
    @("17+33=50", "43-39=4"),
    @("66-37=29", "10+41=51"),
    @("23+42=65", "95-63=32"),
    @("83-24=59", "45+18=63"),
    @("57-11=46", "91-57=34"),
    @("31+34=65", "13+66=79"),
    @("80-42=38", "82-45=37"),
    @("30-29=1", "93-12=81"),
    @("90-34=56", "75-31=44"),
    @("62-38=24", "67+19=86"),
    @("38-16=22", "81-21=60"),
    @("22-0=22", "58+27=85"),
    @("9+32=41", "85+0=85"),
    @("77-41=36", "87-85=2"),
    @("54-35=19", "89-61=28"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $r = $d.Content
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "done"
